# Finished dynamic excel reading: update observed values on Sheet3 and
# refresh the active selections on Sheet1 / Sheet3 left over from editing.

$wb = $excel.ActiveWorkbook

# Sheet3: fill in the two previously-blank "Simulated Outflow" observations;
# every downstream statistic (G2:Q2) is formula-driven and recalculates
# automatically.
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("E2").Value = 0.1
$ws3.Range("E3").Value = 0.12

# Leave the cursor where the author last left it on each sheet.
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$ws1.Range("I1").Select()

$ws3.Activate()
$ws3.Range("F3").Select()
